$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("File")

$values = @(
    "TANG-PHONG-1",
    "TANG-PHONG-2",
    "TANG-PHONG-3",
    "TANG-PHONG-4",
    "TANG-PHONG-5",
    "TANG-PHONG-6",
    "TANG-PHONG-7",
    "TANG-PHONG-8",
    "TANG-PHONG-9",
    "TANG-PHONG-10"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

$ws.Range("E9").Select()
